$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 to 43) holds the "Förändrad" date which was bumped by one day
# (45746 -> 45747, i.e. 2025-03-30 -> 2025-03-31).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45746) {
        $cell.Value = 45747
    }
}
